# add alternative task check, create table for protected documents
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of existing tallies in the "Политика" table (rows 4-15) ---
$ws.Range("F4").Value = 5
$ws.Range("J4").Value = 6

$ws.Range("F5").Value = 4
$ws.Range("H5").Value = 4
$ws.Range("J5").Value = 8

$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1

$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 8

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 6

# --- Widen column C to fit the new table's longer labels ---
# ColumnWidth applies Excel's usual +5px padding, so back the target width (20)
# off by that same offset (5/6 character) to land on exactly width=20 in the XML.
$ws.Columns("C").ColumnWidth = 19.166666666666668

# --- Build the new "Объект защиты" (protected documents) table below the existing one ---

# Header row (row 19): reuse the same formatting as the row-3 header (borders +
# centered/wrapped bold-less header style) by copying it, then overwrite the text.
$ws.Range("C3:G3").Copy()
$ws.Range("C19:G19").PasteSpecial(-4122)
$ws.Range("E19").Value = "Ложные объекты"
$ws.Range("F19").Value = "Отсутствующие объекты"
$ws.Range("G19").Value = "Итого недочетов"

# Data rows (20-31): reuse the bordered data-row formatting from row 4.
$ws.Range("C4:G4").Copy()
$ws.Range("C20:G31").PasteSpecial(-4122)

$data = @(
    @("Объект защиты №4",  "Задание №4",  4, 0, 4),
    @("Объект защиты №5",  "Задание №5",  0, 0, 0),
    @("Объект защиты №6",  "Задание №6",  0, 1, 1),
    @("Объект защиты №7",  "Задание №7",  0, 4, 4),
    @("Объект защиты №8",  "Задание №8",  0, 0, 0),
    @("Объект защиты №9",  "Задание №9",  0, 0, 0),
    @("Объект защиты №10", "Задание №10", 0, 0, 0),
    @("Объект защиты №11", "Задание №11", 0, 0, 0),
    @("Объект защиты №12", "Задание №12", 0, 0, 0),
    @("Объект защиты №13", "Задание №13", 0, 0, 0),
    @("Объект защиты №14", "Задание №14", 0, 2, 2),
    @("Объект защиты №15", "Задание №15", 0, 0, 0)
)

$r = 20
foreach ($row in $data) {
    $ws.Range("C$r").Value = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
    $r = $r + 1
}
